$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 346, pushing all existing data (rows 346-442)
# down to rows 348-444. This mirrors the diff, where the whole data block
# shifted down by two rows and two brand-new rows of weekly data were
# inserted at the top of that block.
$ws.Range("A346:A347").EntireRow.Insert()

# --- New row 346: Naranja / Valencia / Primera, updated week ---
$ws.Range("A346").Value2 = 4
$ws.Range("B346").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C346").Value2 = "Los Lagos"
$ws.Range("D346").Value2 = 44711
$ws.Range("E346").Value2 = 10
$ws.Range("F346").Value2 = "Fruta"
$ws.Range("G346").Value2 = 100102
$ws.Range("H346").Value2 = "Cítricos"
$ws.Range("I346").Value2 = 100102005
$ws.Range("J346").Value2 = "Naranja"
$ws.Range("K346").Value2 = "Valencia"
$ws.Range("L346").Value2 = "Primera"
$ws.Range("M346").Value2 = 300
$ws.Range("N346").Value2 = 16000
$ws.Range("O346").Value2 = 16000
$ws.Range("P346").Value2 = 16000
$ws.Range("Q346").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R346").Value2 = "Región de O'Higgins"
$ws.Range("S346").Value2 = 1067
$ws.Range("T346").Value2 = 15

# --- New row 347: Naranja / Valencia / Segunda, updated week ---
$ws.Range("A347").Value2 = 4
$ws.Range("B347").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C347").Value2 = "Los Lagos"
$ws.Range("D347").Value2 = 44711
$ws.Range("E347").Value2 = 10
$ws.Range("F347").Value2 = "Fruta"
$ws.Range("G347").Value2 = 100102
$ws.Range("H347").Value2 = "Cítricos"
$ws.Range("I347").Value2 = 100102005
$ws.Range("J347").Value2 = "Naranja"
$ws.Range("K347").Value2 = "Valencia"
$ws.Range("L347").Value2 = "Segunda"
$ws.Range("M347").Value2 = 150
$ws.Range("N347").Value2 = 14000
$ws.Range("O347").Value2 = 14000
$ws.Range("P347").Value2 = 14000
$ws.Range("Q347").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R347").Value2 = "Región de O'Higgins"
$ws.Range("S347").Value2 = 933
$ws.Range("T347").Value2 = 15
